$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "GEO_RNASEQ"

$data = New-Object 'object[,]' 18,12
$data[0,1] = "TermSourceRef"
$data[0,2] = "Ontology"
$data[0,3] = "TAN"
$data[0,4] = "Content type (validation)"
$data[0,5] = "Notes during templating"
$data[0,6] = "Target term"
$data[0,7] = "Instruction"
$data[0,8] = "Requirement (m/o/n)"
$data[0,9] = "Value (cv/s/d)"
$data[0,10] = "Additional information"
$data[0,11] = "Review comments"
$data[1,0] = "Source Name"
$data[2,0] = "Sample Name"
$data[3,0] = "Parameter [Library strategy]"
$data[3,1] = "NFDI4PSO:0000035"
$data[3,2] = "NFDI4PSO"
$data[3,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000035"
$data[4,0] = "Parameter [Library Selection]"
$data[4,1] = "NFDI4PSO:0000036"
$data[4,2] = "NFDI4PSO"
$data[4,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000036"
$data[5,0] = "Parameter [Library layout]"
$data[5,1] = "user-specific"
$data[5,2] = "user-specific"
$data[5,3] = "user-specific"
$data[6,0] = "Parameter [Library preparation kit]"
$data[6,1] = "NFDI4PSO:0000037"
$data[6,2] = "NFDI4PSO"
$data[6,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000037"
$data[7,0] = "Parameter [Library preparation kit version]"
$data[7,1] = "NFDI4PSO:0000038"
$data[7,2] = "NFDI4PSO"
$data[7,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000038"
$data[8,0] = "Parameter [Adapter sequence]"
$data[8,1] = "NFDI4PSO:0000039"
$data[8,2] = "NFDI4PSO"
$data[8,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000039"
$data[9,0] = "Parameter [Library RNA amount]"
$data[9,1] = "NFDI4PSO:0000016"
$data[9,2] = "NFDI4PSO"
$data[9,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000016"
$data[10,0] = "Parameter [rRNA depletion]"
$data[10,1] = "NFDI4PSO:0000082"
$data[10,2] = "NFDI4PSO"
$data[10,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000082"
$data[11,0] = "Parameter [Next generation sequencing instrument model]"
$data[11,1] = "NFDI4PSO:0000040"
$data[11,2] = "NFDI4PSO"
$data[11,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000040"
$data[12,0] = "Parameter [Base-calling Software]"
$data[12,1] = "NFDI4PSO:0000017"
$data[12,2] = "NFDI4PSO"
$data[12,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000017"
$data[13,0] = "Parameter [Base-calling Software Version]"
$data[13,1] = "NFDI4PSO:0000018"
$data[13,2] = "NFDI4PSO"
$data[13,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000018"
$data[14,0] = "Parameter [Base-calling Software Parameters]"
$data[14,1] = "NFDI4PSO:0000019"
$data[14,2] = "NFDI4PSO"
$data[14,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000019"
$data[15,0] = "Parameter [Library strand]"
$data[15,1] = "NFDI4PSO:0000020"
$data[15,2] = "NFDI4PSO"
$data[15,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000020"
$data[16,0] = "Parameter [Raw data file format]"
$data[16,1] = "NFDI4PSO:0000021"
$data[16,2] = "NFDI4PSO"
$data[16,3] = "http://purl.obolibrary.org/obo/NFDI4PSO_0000021"
$data[17,0] = "Data File Name"
$ws.Range("A1:L18").Value = $data
$ws.Columns.AutoFit() | Out-Null
$ws.Activate()
